# Updates the crypto price/volume table on Sheet1 with the latest scraped
# values (GitHub Actions refresh run).
#
# Note: a handful of "Price" cells contain values that look like plain
# decimal numbers (e.g. "1.00", "600.20"). Assigning those strings directly
# to Range.Value lets Excel auto-convert them to floating point numbers,
# which silently drops the trailing zeros / exact text ("1.00" -> 1,
# "10.10" -> 10.1). To keep them as literal text (matching how the source
# data is published) we prefix those specific assignments with a leading
# apostrophe, exactly like typing `'1.00` into a cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.491.08'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '3.899.96'
$ws.Range('E3').Value = '  +3.91%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'600.20"
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = "'164.84"
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = '3.897.76'
$ws.Range('E7').Value = '  +3.90%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('E10').Value = '  -3.97%  '
$ws.Range('D11').Value = "'6.38"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = "'36.92"
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('D15').Value = '4.549.06'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').Value = '3.914.86'
$ws.Range('E16').Value = '  +4.28%  '
$ws.Range('D17').Value = '68.656.77'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = "'7.43"
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = "'16.96"
$ws.Range('E20').Value = '  -4.83%  '
$ws.Range('D21').Value = "'11.18"
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('D22').Value = "'485.37"
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('D23').Value = "'0.718"
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('E24').Value = '  +11.11%  '
$ws.Range('D25').Value = "'84.26"
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').Value = "'2.25"
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').Value = "'12.04"
$ws.Range('E27').Value = '  -2.04%  '
$ws.Range('D28').Value = "'10.10"
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('D31').Value = '4.050.44'
$ws.Range('E31').Value = '  +4.00%  '
$ws.Range('D32').Value = "'7.85"
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D33').Value = "'2.38"
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').Value = "'31.91"
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('D35').Value = '3.839.97'
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = "'1.03"
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').Value = "'5.88"
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').Value = "'0.319"
$ws.Range('E41').Value = '  -1.80%  '
$ws.Range('D42').Value = "'438.11"
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('E43').Value = '  -4.61%  '
$ws.Range('D44').Value = "'48.44"
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D47').Value = "'8.46"
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').Value = '2.833.65'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').Value = "'142.10"
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').Value = "'25.91"
$ws.Range('E50').Value = '  +9.62%  '
$ws.Range('E51').Value = '  +0.23%  '
